# The document ends with two consecutive empty paragraphs right before the
# final section break (an extra blank line left over after the closing
# "... provisioning checklist." paragraph). Remove the extra blank line so
# only a single trailing empty paragraph remains, same as Word does when a
# user parks the cursor on the blank line and presses Backspace/Delete.

$d = $word.ActiveDocument

$blank = $null

# Anchor on the closing sentence of the toolkit paragraph so we find the
# right spot even if paragraph indices shift, then step to the very next
# paragraph (the blank one right after it).
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute("provisioning checklist.", $false, $false, $false, $false,
                               $false, $true, 1, $false, "", 0)

if ($found -and $anchor.Find.Found) {
    $anchor.Collapse(0)             # wdCollapseEnd
    $nextPara = $anchor.Next(4, 1)  # wdParagraph, one paragraph forward
    if ($nextPara -ne $null -and ($nextPara.End - $nextPara.Start) -le 1) {
        $blank = $nextPara
    }
}

if ($blank -eq $null) {
    # Fallback: the document ends with two empty paragraphs; drop the first
    # of the two (the one right before the very last paragraph).
    $count = $d.Paragraphs.Count
    $candidate = $d.Paragraphs.Item($count - 1)
    if (($candidate.Range.End - $candidate.Range.Start) -le 1) {
        $blank = $candidate.Range
    }
}

if ($blank -ne $null) {
    $blank.Delete()
}
